# "Add files via upload" — re-upload of the Season 1 standings workbook with
# the placeholder "MasterCard Renault F1 Team" rows (no real driver/team
# assigned, 0 points) removed from both the Drivers' and Constructors'
# standings sheets.
#
# WDC sheet ("SEASON 1 WDC"): two placeholder rows (19th & 20th place) are
# removed; the trailing 3pt-tall spacer row shifts up to become row 20.
#
# WCC sheet ("SEASON 1 WCC"): one placeholder row (10th place) is removed;
# the trailing 6pt-tall spacer row shifts up to become row 11.

$wb = $excel.ActiveWorkbook

$wsDrivers = $wb.Worksheets.Item("SEASON 1 WDC")
$wsDrivers.Rows("20:21").Delete()

$wsConstructors = $wb.Worksheets.Item("SEASON 1 WCC")
$wsConstructors.Rows("11:11").Delete()

# Re-apply the existing sort (data is already in descending-points order) so
# the sheets' persisted sortState/dimension shrink to the new, smaller range
# instead of continuing to reference the now-deleted rows.
$rngDrivers = $wsDrivers.Range("B1:AI20")
$keyDrivers = $wsDrivers.Range("E1:E20")
$wsDrivers.Sort.SortFields.Clear()
$wsDrivers.Sort.SortFields.Add($keyDrivers, 0, 2, 0, 0) | Out-Null
$wsDrivers.Sort.SetRange($rngDrivers)
$wsDrivers.Sort.Header = 0
$wsDrivers.Sort.Apply()

$rngConstructors = $wsConstructors.Range("B2:AH11")
$keyConstructors = $wsConstructors.Range("AC2:AC11")
$wsConstructors.Sort.SortFields.Clear()
$wsConstructors.Sort.SortFields.Add($keyConstructors, 0, 2, 0, 0) | Out-Null
$wsConstructors.Sort.SetRange($rngConstructors)
$wsConstructors.Sort.Header = 0
$wsConstructors.Sort.Apply()

# Re-write the constructors' point-total formulas as one range assignment so
# Excel collapses them back into a single shared formula (as they were
# before the extra placeholder row existed).
$wsConstructors.Range("AC2:AC10").Formula = "=SUM(D2:AB2)"
